$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts CardCode.. and all other columns one to the right)
$ws.Columns.Item(1).EntireColumn.Insert()

# New header in the freshly inserted column A
$ws.Range("A1").Value = "StaffId"

# Match the author's narrower width for the new StaffId column (stored width = 8)
$ws.Columns.Item(1).ColumnWidth = 7.1666666666667

# Update the view: clear the old frozen/top-left scroll position and select D7
$ws.Range("D7").Select()
